$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 115; this shifts existing rows 115-119 down to 116-120,
# preserving all their data/formatting (including the D column's date style).
$ws.Rows.Item(115).Insert()

# Populate the newly inserted row 115 with the new weekly record.
$ws.Cells.Item(115, 1).Value = 4
$ws.Cells.Item(115, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(115, 3).Value = "Los Lagos"
$ws.Cells.Item(115, 4).Value = 44509
$ws.Cells.Item(115, 4).NumberFormat = $ws.Cells.Item(116, 4).NumberFormat
$ws.Cells.Item(115, 5).Value = 10
$ws.Cells.Item(115, 6).Value = 100112009
$ws.Cells.Item(115, 7).Value = "Acelga"
$ws.Cells.Item(115, 8).Value = "Sin especificar"
$ws.Cells.Item(115, 9).Value = "Primera"
$ws.Cells.Item(115, 10).Value = 200
$ws.Cells.Item(115, 11).Value = 3000
$ws.Cells.Item(115, 12).Value = 3000
$ws.Cells.Item(115, 13).Value = 3000
$ws.Cells.Item(115, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(115, 15).Value = "Región del Maule"
$ws.Cells.Item(115, 16).Value = 750
$ws.Cells.Item(115, 17).Value = 4
$ws.Cells.Item(115, 18).Value = "Hortaliza"
